$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Update existing row 2 values (lot LOTE1063 -> LOTE1095, with associated date/remito/pedido updates)
Set-TextValue "I2" "LOTE1095"
Set-TextValue "J2" "20200107"
Set-TextValue "K2" "20210829"
Set-TextValue "L2" "1095"
$ws.Range("M2").Value = 18
Set-TextValue "N2" "0000-00001095"
Set-TextValue "Y2" "20191108"
Set-TextValue "Z2" "3406778"
$ws.Range("AD2").Value = 18

# Column I (out_lote) widened slightly to fit the new best-fit content
$ws.Columns("I:I").ColumnWidth = 8.8333333333

# Add new row 3 with replicated values
$ws.Range("A3").Value = "Calidad"
$ws.Range("B3").Value = 1002
$ws.Range("E3").Value = "ON"
$ws.Range("M3").Value = 19
$ws.Range("O3").Value = "OFF"
$ws.Range("AD3").Value = 19
$ws.Range("AG3").Value = 2032500

# Update view: set selection to A2 (also resets the scrolled top-left cell to A1)
$ws.Range("A2").Select()
